$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new blank row at row 162 (shifts rows 162..177 down to 163..178)
$ws.Rows.Item(162).Insert()

# 2. Copy the formatting of row 160 (A160:Z160) onto the new row 162 so it
#    picks up the same style set used by the "55900xxx" skill rows.
$ws.Range("A160:Z160").Copy()
$ws.Range("A162:Z162").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Row height for the new row matches the surrounding "55990xxx" rows (36pt)
$ws.Rows.Item(162).RowHeight = 36

# 4. Populate the new row 162 ("piaoqie" / 剽窃 - card thief skill)
$ws.Range("A162").Value2 = 55900035
$ws.Range("B162").Value2 = "剽窃"
$ws.Range("C162").Value2 = "特殊"
$ws.Range("H162").Value2 = "s.Owner.AddRandomCardJob(s.Rival.Job,s.Level);"
$ws.Range("Q162").Value2 = "Active"
$ws.Range("R162").Value2 = "true"
$ws.Range("S162").Value2 = "召唤时获得一张随机对方职业卡牌"
$ws.Range("X162").Value2 = 14
$ws.Range("Y162").Value2 = "piaoqie"

# 5. Resize the worksheet table ("表3_25") so it covers the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:Z178"))

# 6. Update the view selection to match the edited area
$ws.Range("F161").Select()

Write-Output "done"
